$d = $word.ActiveDocument

# --- Body (document.xml) replacements ---
# " 000110364901 - 2 " -> "  "
$d.Content.Find.Execute(" 000110364901 - 2 ", $true, $false, $false, $false, $false, $true, 1, $false, "  ", 2)

# "QWR" (bold run, body) -> "TERE"
$d.Content.Find.Execute("QWR", $true, $false, $false, $false, $false, $true, 1, $false, "TERE", 2)

# --- Header (header1.xml) replacements ---
$header = $d.Sections(1).Headers(1).Range

# "QWER" -> "TRE"
$header.Find.Execute("QWER", $true, $false, $false, $false, $false, $true, 1, $false, "TRE", 2)

# "QWR" -> "TERE"
$header.Find.Execute("QWR", $true, $false, $false, $false, $false, $true, 1, $false, "TERE", 2)

# "Qwer" -> "Tre" (5 occurrences) -- Replace All
$header.Find.Execute("Qwer", $true, $false, $false, $false, $false, $true, 1, $false, "Tre", 2)

# "qwer" -> "tre" (3 occurrences) -- Replace All
$header.Find.Execute("qwer", $true, $false, $false, $false, $false, $true, 1, $false, "tre", 2)
